$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'322.99"
$ws.Range("E2").Value = "'-1.87%"
$ws.Range("D3").Value = "'39.72"
$ws.Range("E3").Value = "'-0.76%"
$ws.Range("D4").Value = "'5.875"
$ws.Range("E4").Value = "'11.36%"
$ws.Range("D5").Value = "'0.08029"
$ws.Range("E5").Value = "'-1.02%"
$ws.Range("D6").Value = "'8.656"
$ws.Range("E7").Value = "'1.31%"
$ws.Range("D8").Value = "'0.9301"
$ws.Range("E8").Value = "'-0.48%"
$ws.Range("D9").Value = "'0.1225"
$ws.Range("E9").Value = "'-8.66%"
$ws.Range("D10").Value = "'0.1957"
$ws.Range("E10").Value = "'0.08%"
$ws.Range("D11").Value = "'8.752"
$ws.Range("E11").Value = "'20.28%"
$ws.Range("D12").Value = "'0.09111"
$ws.Range("E12").Value = "'-0.74%"
$ws.Range("D13").Value = "'0.03531"
$ws.Range("E13").Value = "'2.67%"
$ws.Range("D14").Value = "'0.09563"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("D15").Value = "'0.001293"
$ws.Range("E15").Value = "'-7.03%"
$ws.Range("D16").Value = "'0.006177"
$ws.Range("E16").Value = "'2.72%"
$ws.Range("D17").Value = "'3.352"
$ws.Range("E17").Value = "'-0.25%"
$ws.Range("D18").Value = "'4.571"
$ws.Range("E18").Value = "'0.96%"
$ws.Range("D19").Value = "'2.951"
$ws.Range("E19").Value = "'-0.27%"
$ws.Range("E20").Value = "'0.94%"
$ws.Range("E21").Value = "'7.99%"
$ws.Range("E22").Value = "'4.20%"
$ws.Range("D23").Value = "'0.04406"
$ws.Range("E23").Value = "'-0.82%"
$ws.Range("E24").Value = "'3.15%"
$ws.Range("D25").Value = "'0.004389"
$ws.Range("E25").Value = "'0.66%"
$ws.Range("E26").Value = "'-11.64%"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("D39").Value = "'0.02423"
$ws.Range("E39").Value = "'-3.46%"
$ws.Range("D40").Value = "'0.05224"
$ws.Range("E40").Value = "'-0.60%"
$ws.Range("D41").Value = "'0.007447"
$ws.Range("E41").Value = "'-3.46%"
$ws.Range("D42").Value = "'0.009371"
$ws.Range("E42").Value = "'10.09%"
$ws.Range("E43").Value = "'-1.80%"
$ws.Range("D44").Value = "'0.002120"
$ws.Range("D45").Value = "'0.01125"
$ws.Range("E45").Value = "'38.32%"
$ws.Range("D46").Value = "'0.00006722"
$ws.Range("E46").Value = "'0.80%"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E48").Value = "'5.25%"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E51").Value = "'-0.01%"
